$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Author name
$ws.Range("A2").Value = "Adrien Allemand"

# Insert a new blank row before the old "Total" row (row 32) so the
# total row moves from row 32 to row 33, and copy the formatting from
# the row above (row 31) into the freshly inserted row 32.
$ws.Rows.Item(32).Insert()
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fix up the SUM formula so it covers the new data range.
$ws.Range("C33").Formula = "=SUM(C5:C32)"

# Update the journal entries (A column values are Excel date serials;
# cells already carry the dd/mm/yy number format).
$ws.Range("A5").Value = 43149
$ws.Range("B5").Value = "Discution de Groupe du projet"
$ws.Range("C5").Value = 3
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = 43150
$ws.Range("B6").Value = "Relecture de la synthèse du brainstorming initial"
$ws.Range("C6").Value = 0.5
$ws.Rows.Item(6).AutoFit()

$ws.Range("A7").Value = 43157
$ws.Range("B7").Value = "Séance de retour sur la proposition du projet"
$ws.Range("C7").Value = 0.5

$ws.Range("A8").Value = 43157
$ws.Range("B8").Value = "Organisation interne suite a la séance de retour"
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = 43158
$ws.Range("B9").Value = "spécification du projet, création d'un shéma de dépendence des fonctionalitées"
$ws.Range("C9").Value = 5
$ws.Rows.Item(9).RowHeight = 30

$ws.Range("A10").Value = 43163
$ws.Range("B10").Value = "Call pour le cahier des charges"
$ws.Range("C10").Value = 1

# Match the cursor position left behind by the editing session.
[void]$ws.Range("F14").Select()
